$d = $word.ActiveDocument

# "and says" -> "and says," before the quoted dialogue
$d.Content.Find.Execute("and says", $true, $false, $false, $false, $false, $true, 1, $false, "and says,", 2)

# Overlord 3 puzzle description: "maze type" -> "laser puzzles"
$d.Content.Find.Execute("maze type", $true, $false, $false, $false, $false, $true, 1, $false, "laser puzzles", 2)

# Rename the character "boy" -> "kid" throughout the story
$d.Content.Find.Execute("boy", $true, $false, $false, $false, $false, $true, 1, $false, "kid", 2)
